# Edit: "Load distribution" sheet — expand the load-distribution table so
# every node index (1-24) has a row. Previously nodes 11, 12, 17, 21, 22,
# 23 and 24 had no entry; this inserts placeholder rows for them
# (Load # = 0, % of system load = 0) so "all indexes" line up and "all
# variables [are] created for all nodes", while leaving the untouched
# existing rows exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Load distribution")

# Node 10 (row 11) is followed by node 13 (row 12) in the original sheet
# -- insert two blank rows above row 12 to hold the missing nodes 11 & 12.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12, 1).Value = 0
$ws.Cells.Item(12, 2).Value = 11
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = 0
$ws.Cells.Item(13, 2).Value = 12
$ws.Cells.Item(13, 3).Value = 0

# Node 16 (now row 17) is followed by node 18 (now row 18) -- insert a
# blank row above it to hold the missing node 17.
$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18, 1).Value = 0
$ws.Cells.Item(18, 2).Value = 17
$ws.Cells.Item(18, 3).Value = 0

# Node 20 is the last entry (now row 21) -- append placeholder rows for
# the missing nodes 21, 22, 23 and 24.
$ws.Cells.Item(22, 1).Value = 0
$ws.Cells.Item(22, 2).Value = 21
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(23, 1).Value = 0
$ws.Cells.Item(23, 2).Value = 22
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(24, 1).Value = 0
$ws.Cells.Item(24, 2).Value = 23
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(25, 1).Value = 0
$ws.Cells.Item(25, 2).Value = 24
$ws.Cells.Item(25, 3).Value = 0

# Activate the "Load distribution" sheet and leave the selection on the
# first empty row below the table (C26), matching where the cursor ends
# up after typing the last new row of data.
$ws.Activate()
$ws.Range("C26").Select()
